$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the unnecessary "PresenationXMLChanges.sql" bullet (its SQL
#    changes were folded directly into InstallData.xml per the commit msg).
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "PresenationXMLChanges\.sql") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Relocate the stray "_GoBack" bookmark (left over from the author's last
#    edit position) so it sits at the start of the "8. Edit the ORNG items"
#    heading instead of after "...instead of 80" near the end of the doc.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^8\.") {
        $ip = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $ip)
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Fix the accidentally-split sentence about telnet-testing the local
#    machine name: the words "name for your server..." were duplicated
#    across a paragraph break. Re-join them into one flowing sentence that
#    ends with "...to that name and port:" followed by the existing
#    "C:>telnet ..." paragraph.
# ---------------------------------------------------------------------------
$oldWhole = " you want to use the local machine " + `
    "name for your server, this might be something different than what shows up in the URL that you use. The way you can test is by first starting up Tomcat and then attempting to telnet " + `
    "to that name and port:"

$rFull = $d.Content
$rFull.Find.ClearFormatting()
$found = $rFull.Find.Execute($oldWhole, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rFull.Start
    $rFull.Delete()

    $t1 = " you want to use the local machine name for your server, this might be something different than what shows up in the URL that you use. The way you can test is by first starting up Tomcat and then attempting to telnet "
    $t2 = "to that name and port:"

    $ip1 = $d.Range($start, $start)
    $ip1.InsertAfter($t1)

    $ip2 = $d.Range($start + $t1.Length, $start + $t1.Length)
    $ip2.InsertAfter($t2)
}
